$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2026-02-21 Saturday" "2026-02-22 Sunday"

Replace-Text "43÷7=6, 1" "39÷9=4, 3"
Replace-Text "38÷2=19, 0" "46÷2=23, 0"
Replace-Text "13÷9=1, 4" "31÷5=6, 1"
Replace-Text "53÷8=6, 5" "88÷7=12, 4"
Replace-Text "64÷8=8, 0" "25÷3=8, 1"

Replace-Text "59÷7=8, 3" "71÷4=17, 3"
Replace-Text "28÷9=3, 1" "26÷7=3, 5"
Replace-Text "65÷9=7, 2" "16÷3=5, 1"
Replace-Text "20÷7=2, 6" "21÷9=2, 3"
Replace-Text "11÷6=1, 5" "35÷2=17, 1"

Replace-Text "63÷5=12, 3" "52÷4=13, 0"
Replace-Text "65÷3=21, 2" "64÷9=7, 1"
Replace-Text "38÷9=4, 2" "44÷2=22, 0"
Replace-Text "58÷7=8, 2" "74÷7=10, 4"
Replace-Text "68÷2=34, 0" "64÷9=7, 1"

Replace-Text "95÷7=13, 4" "16÷5=3, 1"
Replace-Text "25÷5=5, 0" "34÷8=4, 2"
Replace-Text "94÷8=11, 6" "32÷4=8, 0"
Replace-Text "81÷8=10, 1" "84÷3=28, 0"
Replace-Text "76÷9=8, 4" "36÷5=7, 1"

Replace-Text "41÷4=10, 1" "22÷7=3, 1"
Replace-Text "24÷7=3, 3" "63÷4=15, 3"
Replace-Text "30÷7=4, 2" "64÷2=32, 0"
Replace-Text "52÷2=26, 0" "92÷2=46, 0"
Replace-Text "23÷2=11, 1" "94÷9=10, 4"

Write-Output "Done"
